$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp6"
$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 23.061728
$ws.Range("H2").Value = 69.18518399999999
$ws.Range("I2").Value = 0.6130144106248721
$ws.Range("J2").Value = 0.6902769593117909
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.470843666666667
$ws.Range("N2").Value = 10.412531
$ws.Range("O2").Value = 0.05653984104486641
$ws.Range("P2").Value = 0.0651263661336549
$ws.Range("Q2").Value = 80.04365257118934
$ws.Range("R2").Value = 720.3928731407041
$ws.Range("S2").Value = 0.03465973733494273
$ws.Range("T2").Value = 0.0449552299857657

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp6"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 23.061728
$ws.Range("H3").Value = 69.18518399999999
$ws.Range("I3").Value = 0.6130144106248721
$ws.Range("J3").Value = 0.6902769593117909
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 32.24261766666667
$ws.Range("N3").Value = 96.72785300000001
$ws.Range("O3").Value = 0.5252303626496961
$ws.Range("P3").Value = 0.6049954204026234
$ws.Range("Q3").Value = 743.5704786366613
$ws.Range("R3").Value = 6692.134307729952
$ws.Range("S3").Value = 0.3219737812019913
$ws.Range("T3").Value = 0.4176143991930815

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp6"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 23.061728
$ws.Range("H4").Value = 69.18518399999999
$ws.Range("I4").Value = 0.6130144106248721
$ws.Range("J4").Value = 0.6902769593117909
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9959919999999999
$ws.Range("N4").Value = 2.987976
$ws.Range("O4").Value = 0.01622465163233374
$ws.Range("P4").Value = 0.01868863765923708
$ws.Range("Q4").Value = 22.969296594176
$ws.Range("R4").Value = 206.723669347584
$ws.Range("S4").Value = 0.009945945257988939
$ws.Range("T4").Value = 0.012900335977098

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bmp6"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.061728
$ws.Range("H5").Value = 69.18518399999999
$ws.Range("I5").Value = 0.6130144106248721
$ws.Range("J5").Value = 0.6902769593117909
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3973579999999999
$ws.Range("N5").Value = 1.192074
$ws.Range("O5").Value = 0.006472938661476068
$ws.Range("P5").Value = 0.00745596318343835
$ws.Range("Q5").Value = 9.163762114623998
$ws.Range("R5").Value = 82.47385903161599
$ws.Range("S5").Value = 0.0039680046785757
$ws.Range("T5").Value = 0.005146679595004484

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Bmp6"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.061728
$ws.Range("H6").Value = 69.18518399999999
$ws.Range("I6").Value = 0.6130144106248721
$ws.Range("J6").Value = 0.6902769593117909
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.2807625
$ws.Range("N6").Value = 48.561525
$ws.Range("O6").Value = 0.3955322060116276
$ws.Range("P6").Value = 0.3037336126210463
$ws.Range("Q6").Value = 559.9563404076
$ws.Range("R6").Value = 3359.7380424456
$ws.Range("S6").Value = 0.2424669421513734
$ws.Range("T6").Value = 0.2096603145608412

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp6"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.926013333333334
$ws.Range("H7").Value = 5.778040000000001
$ws.Range("I7").Value = 0.05119624723650278
$ws.Range("J7").Value = 0.05764887294340218
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.470843666666667
$ws.Range("N7").Value = 10.412531
$ws.Range("O7").Value = 0.05653984104486641
$ws.Range("P7").Value = 0.0651263661336549
$ws.Range("Q7").Value = 6.684891179915557
$ws.Range("R7").Value = 60.16402061924001
$ws.Range("S7").Value = 0.002894627680845548
$ws.Range("T7").Value = 0.003754461606504562

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Bmp6"
$ws.Range("C8").Value = "Bmpr1a"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.926013333333334
$ws.Range("H8").Value = 5.778040000000001
$ws.Range("I8").Value = 0.05119624723650278
$ws.Range("J8").Value = 0.05764887294340218
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 32.24261766666667
$ws.Range("N8").Value = 96.72785300000001
$ws.Range("O8").Value = 0.5252303626496961
$ws.Range("P8").Value = 0.6049954204026234
$ws.Range("Q8").Value = 62.0997115275689
$ws.Range("R8").Value = 558.8974037481202
$ws.Range("S8").Value = 0.02688982350233186
$ws.Range("T8").Value = 0.03487730412213103

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Bmp6"
$ws.Range("C9").Value = "Bmpr1a"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.926013333333334
$ws.Range("H9").Value = 5.778040000000001
$ws.Range("I9").Value = 0.05119624723650278
$ws.Range("J9").Value = 0.05764887294340218
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.9959919999999999
$ws.Range("N9").Value = 2.987976
$ws.Range("O9").Value = 0.01622465163233374
$ws.Range("P9").Value = 0.01868863765923708
$ws.Range("Q9").Value = 1.918293871893333
$ws.Range("R9").Value = 17.26464484704
$ws.Range("S9").Value = 0.0008306412762950867
$ws.Range("T9").Value = 0.00107737889790264

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Bmp6"
$ws.Range("C10").Value = "Bmpr1a"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.926013333333334
$ws.Range("H10").Value = 5.778040000000001
$ws.Range("I10").Value = 0.05119624723650278
$ws.Range("J10").Value = 0.05764887294340218
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3973579999999999
$ws.Range("N10").Value = 1.192074
$ws.Range("O10").Value = 0.006472938661476068
$ws.Range("P10").Value = 0.00745596318343835
$ws.Range("Q10").Value = 0.7653168061066666
$ws.Range("R10").Value = 6.88785125496
$ws.Range("S10").Value = 0.0003313901680596461
$ws.Range("T10").Value = 0.0004298278742327219

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Bmp6"
$ws.Range("C11").Value = "Bmpr1a"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.926013333333334
$ws.Range("H11").Value = 5.778040000000001
$ws.Range("I11").Value = 0.05119624723650278
$ws.Range("J11").Value = 0.05764887294340218
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 24.2807625
$ws.Range("N11").Value = 48.561525
$ws.Range("O11").Value = 0.3955322060116276
$ws.Range("P11").Value = 0.3037336126210463
$ws.Range("Q11").Value = 46.76507231850001
$ws.Range("R11").Value = 280.590433911
$ws.Range("S11").Value = 0.02024976460897064
$ws.Range("T11").Value = 0.01750990044263123

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Bmp6"
$ws.Range("C12").Value = "Bmpr1a"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 12.632464
$ws.Range("H12").Value = 25.264928
$ws.Range("I12").Value = 0.3357893421386252
$ws.Range("J12").Value = 0.2520741677448068
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.470843666666667
$ws.Range("N12").Value = 10.412531
$ws.Range("O12").Value = 0.05653984104486641
$ws.Range("P12").Value = 0.0651263661336549
$ws.Range("Q12").Value = 43.84530766879467
$ws.Range("R12").Value = 263.071846012768
$ws.Range("S12").Value = 0.01898547602907813
$ws.Range("T12").Value = 0.01641667454138463

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Bmp6"
$ws.Range("C13").Value = "Bmpr1a"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 12.632464
$ws.Range("H13").Value = 25.264928
$ws.Range("I13").Value = 0.3357893421386252
$ws.Range("J13").Value = 0.2520741677448068
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 32.24261766666667
$ws.Range("N13").Value = 96.72785300000001
$ws.Range("O13").Value = 0.5252303626496961
$ws.Range("P13").Value = 0.6049954204026234
$ws.Range("Q13").Value = 407.3037069399306
$ws.Range("R13").Value = 2443.822241639584
$ws.Range("S13").Value = 0.176366757945373
$ws.Range("T13").Value = 0.1525037170874108

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Bmp6"
$ws.Range("C14").Value = "Bmpr1a"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 12.632464
$ws.Range("H14").Value = 25.264928
$ws.Range("I14").Value = 0.3357893421386252
$ws.Range("J14").Value = 0.2520741677448068
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.9959919999999999
$ws.Range("N14").Value = 2.987976
$ws.Range("O14").Value = 0.01622465163233374
$ws.Range("P14").Value = 0.01868863765923708
$ws.Range("Q14").Value = 12.581833084288
$ws.Range("R14").Value = 75.49099850572799
$ws.Range("S14").Value = 0.00544806509804972
$ws.Range("T14").Value = 0.004710922784236442

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Bmp6"
$ws.Range("C15").Value = "Bmpr1a"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 12.632464
$ws.Range("H15").Value = 25.264928
$ws.Range("I15").Value = 0.3357893421386252
$ws.Range("J15").Value = 0.2520741677448068
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.3973579999999999
$ws.Range("N15").Value = 1.192074
$ws.Range("O15").Value = 0.006472938661476068
$ws.Range("P15").Value = 0.00745596318343835
$ws.Range("Q15").Value = 5.019610630111998
$ws.Range("R15").Value = 30.11766378067199
$ws.Range("S15").Value = 0.002173543814840722
$ws.Range("T15").Value = 0.001879455714201143

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Bmp6"
$ws.Range("C16").Value = "Bmpr1a"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 12.632464
$ws.Range("H16").Value = 25.264928
$ws.Range("I16").Value = 0.3357893421386252
$ws.Range("J16").Value = 0.2520741677448068
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 24.2807625
$ws.Range("N16").Value = 48.561525
$ws.Range("O16").Value = 0.3955322060116276
$ws.Range("P16").Value = 0.3037336126210463
$ws.Range("Q16").Value = 306.7258581738
$ws.Range("R16").Value = 1226.9034326952
$ws.Range("S16").Value = 0.1328154992512836
$ws.Range("T16").Value = 0.07656339761757379

